$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: TOTAL (5-17 ans) - value changes only
$ws.Cells.Item(2, 3).Value = 1212250

# Row 3: non_pdi -> hote
$ws.Cells.Item(3, 1).Value = "hote (5-17 y.o.)"
$ws.Cells.Item(3, 2).Value = "hote"
$ws.Cells.Item(3, 3).Value = 1079110

# Row 4: pdi -> idp_host
$ws.Cells.Item(4, 1).Value = "idp_host (5-17 y.o.)"
$ws.Cells.Item(4, 2).Value = "idp_host"
$ws.Cells.Item(4, 3).Value = 66749

# Row 5: Filles (5-17 ans) -> retourne (5-17 y.o.)
$ws.Cells.Item(5, 1).Value = "retourne (5-17 y.o.)"
$ws.Cells.Item(5, 2).Value = "retourne"
$ws.Cells.Item(5, 3).Value = 49422

# Row 6: Garcons (5-17 ans) -> idp_site (5-17 y.o.)
$ws.Cells.Item(6, 1).Value = "idp_site (5-17 y.o.)"
$ws.Cells.Item(6, 2).Value = "idp_site"
$ws.Cells.Item(6, 3).Value = 16970

# Row 7: Éducation préscolaire (5 ans) -> Filles (5-17 ans)
$ws.Cells.Item(7, 1).Value = "Filles (5-17 ans)"
$ws.Cells.Item(7, 2).Value = "Tous les groupes de population"
$ws.Cells.Item(7, 3).Value = 618247

# Row 8: new row - Garcons (5-17 ans)
$ws.Cells.Item(8, 1).Value = "Garcons (5-17 ans)"
$ws.Cells.Item(8, 2).Value = "Tous les groupes de population"
$ws.Cells.Item(8, 3).Value = 594002

# Row 9: new row - Éducation préscolaire (5 ans)
$ws.Cells.Item(9, 1).Value = "Éducation préscolaire (5 ans)"
$ws.Cells.Item(9, 2).Value = "Tous les groupes de population"
$ws.Cells.Item(9, 3).Value = 71566

# Row 10: Enfants en situation de handicap (shifted from row 8)
$ws.Cells.Item(10, 1).Value = "Enfants en situation de handicap"
$ws.Cells.Item(10, 2).Value = "Tous les groupes de population"
$ws.Cells.Item(10, 3).Value = 121225
